$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cells used to copy cell formatting (style) for brand-new rows
$templateA = $ws.Cells.Item(2, 1)
$templateB = $ws.Cells.Item(2, 2)
$templateC = $ws.Cells.Item(2, 3)

# Row 2
$ws.Cells.Item(2, 3).Value = 8515.0

# Row 3
$ws.Cells.Item(3, 3).Value = 1061.0

# Row 4
$ws.Cells.Item(4, 3).Value = 767.0

# Row 5
$ws.Cells.Item(5, 3).Value = 518.0

# Row 6
$ws.Cells.Item(6, 3).Value = 458.0

# Row 7
$ws.Cells.Item(7, 3).Value = 411.0

# Row 9
$ws.Cells.Item(9, 3).Value = 390.0

# Row 11
$ws.Cells.Item(11, 3).Value = 332.0

# Row 12
$ws.Cells.Item(12, 3).Value = 327.0

# Row 13
$ws.Cells.Item(13, 3).Value = 236.0

# Row 14
$ws.Cells.Item(14, 3).Value = 183.0

# Row 15
$ws.Cells.Item(15, 3).Value = 181.0

# Row 16
$ws.Cells.Item(16, 3).Value = 158.0

# Row 17
$ws.Cells.Item(17, 3).Value = 149.0

# Row 19
$ws.Cells.Item(19, 3).Value = 139.0

# Row 20
$ws.Cells.Item(20, 3).Value = 131.0

# Row 23
$ws.Cells.Item(23, 3).Value = 104.0

# Row 25
$ws.Cells.Item(25, 3).Value = 86.0

# Row 26
$ws.Cells.Item(26, 3).Value = 79.0

# Row 27
$ws.Cells.Item(27, 3).Value = 73.0

# Row 35
$ws.Cells.Item(35, 1).Value = 'Regione Autonoma Valle D''Aosta'
$ws.Cells.Item(35, 2).NumberFormat = "@"
$ws.Cells.Item(35, 2).Value = '80002270074'
$ws.Cells.Item(35, 3).Value = 42.0

# Row 36
$ws.Cells.Item(36, 1).Value = 'Consorzio I.T. Srl'
$ws.Cells.Item(36, 2).NumberFormat = "@"
$ws.Cells.Item(36, 2).Value = '01321400192'
$ws.Cells.Item(36, 3).Value = 40.0

# Row 37
$ws.Cells.Item(37, 1).Value = 'ROMA CAPITALE'
$ws.Cells.Item(37, 2).NumberFormat = "@"
$ws.Cells.Item(37, 2).Value = '02438750586'
$ws.Cells.Item(37, 3).Value = 39.0

# Row 38
$ws.Cells.Item(38, 1).Value = 'UNIMATICA S.P.A'
$ws.Cells.Item(38, 2).NumberFormat = "@"
$ws.Cells.Item(38, 2).Value = '02098391200'
$ws.Cells.Item(38, 3).Value = 38.0

# Row 39
$ws.Cells.Item(39, 1).Value = 'SI.net Servizi Informatici S.r.L.'
$ws.Cells.Item(39, 2).NumberFormat = "@"
$ws.Cells.Item(39, 2).Value = '02743730125'
$ws.Cells.Item(39, 3).Value = 36.0

# Row 40
$ws.Cells.Item(40, 1).Value = 'Unicredit, Societa'' per Azioni'
$ws.Cells.Item(40, 2).NumberFormat = "@"
$ws.Cells.Item(40, 2).Value = '00348170101'
$ws.Cells.Item(40, 3).Value = 33.0

# Row 41
$ws.Cells.Item(41, 3).Value = 32.0

# Row 42
$ws.Cells.Item(42, 1).Value = 'DCS SOFTWARE E SERVIZI S.R.L.'
$ws.Cells.Item(42, 2).NumberFormat = "@"
$ws.Cells.Item(42, 2).Value = '08063140019'
$ws.Cells.Item(42, 3).Value = 28.0

# Row 43
$ws.Cells.Item(43, 1).Value = 'Regione Liguria'
$ws.Cells.Item(43, 2).NumberFormat = "@"
$ws.Cells.Item(43, 2).Value = '00849050109'
$ws.Cells.Item(43, 3).Value = 25.0

# Row 44
$ws.Cells.Item(44, 1).Value = 'Novares Spa'
$ws.Cells.Item(44, 2).NumberFormat = "@"
$ws.Cells.Item(44, 2).Value = '12105121003'
$ws.Cells.Item(44, 3).Value = 20.0

# Row 45
$ws.Cells.Item(45, 1).Value = 'Citta'' Metropolitana di Roma Capitale'
$ws.Cells.Item(45, 2).NumberFormat = "@"
$ws.Cells.Item(45, 2).Value = '80034390585'

# Row 48
$ws.Cells.Item(48, 1).Value = 'Nexi SpA'
$ws.Cells.Item(48, 2).NumberFormat = "@"
$ws.Cells.Item(48, 2).Value = '13212880150'
$ws.Cells.Item(48, 3).Value = 18.0

# Row 50
$ws.Cells.Item(50, 1).Value = 'Comune di Palermo'
$ws.Cells.Item(50, 2).NumberFormat = "@"
$ws.Cells.Item(50, 2).Value = '80016350821'
$ws.Cells.Item(50, 3).Value = 17.0

# Row 51
$ws.Cells.Item(51, 1).Value = 'Servizi Locali SpA'
$ws.Cells.Item(51, 2).NumberFormat = "@"
$ws.Cells.Item(51, 2).Value = '03170580751'
$ws.Cells.Item(51, 3).Value = 14.0

# Row 52
$ws.Cells.Item(52, 1).Value = 'Crédit Agricole Group Solutions Società Consortile per azioni'
$ws.Cells.Item(52, 2).NumberFormat = "@"
$ws.Cells.Item(52, 2).Value = '02771790348'
$ws.Cells.Item(52, 3).Value = 12.0

# Row 53
$ws.Cells.Item(53, 1).Value = 'Si.Form Consulting srl'
$ws.Cells.Item(53, 2).NumberFormat = "@"
$ws.Cells.Item(53, 2).Value = '03943960827'

# Row 56
$ws.Cells.Item(56, 3).Value = 7.0

# Row 57
$ws.Cells.Item(57, 1).Value = 'ARCA Servizi s.r.l'
$ws.Cells.Item(57, 2).NumberFormat = "@"
$ws.Cells.Item(57, 2).Value = '09106071005'
$ws.Cells.Item(57, 3).Value = 7.0

# Row 58
$ws.Cells.Item(58, 1).Value = 'Comune di Catania'
$ws.Cells.Item(58, 2).NumberFormat = "@"
$ws.Cells.Item(58, 2).Value = '00137020871'
$ws.Cells.Item(58, 3).Value = 6.0

# Row 59
$ws.Cells.Item(59, 1).Value = 'e-SED Società Cooperativa'
$ws.Cells.Item(59, 2).NumberFormat = "@"
$ws.Cells.Item(59, 2).Value = '02695640421'

# Row 60
$ws.Cells.Item(60, 1).Value = 'ISWEB S.p.A.'
$ws.Cells.Item(60, 2).NumberFormat = "@"
$ws.Cells.Item(60, 2).Value = '01722270665'

# Row 61
$ws.Cells.Item(61, 1).Value = 'Linea Comune Spa'
$ws.Cells.Item(61, 2).NumberFormat = "@"
$ws.Cells.Item(61, 2).Value = '05591950489'

# Row 62
$ws.Cells.Item(62, 1).Value = 'Phoenix IT Solutions S.r.L'
$ws.Cells.Item(62, 2).NumberFormat = "@"
$ws.Cells.Item(62, 2).Value = '07623321218'

# Row 63
$ws.Cells.Item(63, 1).Value = 'CityPoste Payment Digital S.r.l.'
$ws.Cells.Item(63, 2).NumberFormat = "@"
$ws.Cells.Item(63, 2).Value = '02003750672'

# Row 64
$ws.Cells.Item(64, 1).Value = 'ICCREA Banca SpA'
$ws.Cells.Item(64, 2).NumberFormat = "@"
$ws.Cells.Item(64, 2).Value = '04774801007'
$ws.Cells.Item(64, 3).Value = 2.0

# Row 65
$ws.Cells.Item(65, 1).Value = 'Ministero dello Sviluppo Economico'
$ws.Cells.Item(65, 2).NumberFormat = "@"
$ws.Cells.Item(65, 2).Value = '80230390587'

# Row 66
$ws.Cells.Item(66, 1).Value = 'BANCA MONTE DEI PASCHI DI SIENA'
$ws.Cells.Item(66, 2).NumberFormat = "@"
$ws.Cells.Item(66, 2).Value = '00884060526'

# Row 67
$ws.Cells.Item(67, 1).Value = 'Engineering Ingegneria Informatica SpA'
$ws.Cells.Item(67, 2).NumberFormat = "@"
$ws.Cells.Item(67, 2).Value = '00967720285'

# Row 68
$ws.Cells.Item(68, 1).Value = 'Argentea S.r.l.'
$ws.Cells.Item(68, 2).NumberFormat = "@"
$ws.Cells.Item(68, 2).Value = '02260390220'

# Row 69
$ws.Cells.Item(69, 1).Value = 'Società Almaviva S.p.A.'
$ws.Cells.Item(69, 2).NumberFormat = "@"
$ws.Cells.Item(69, 2).Value = '08450891000'

# Row 70
$ws.Cells.Item(70, 1).Value = 'I.C.A. - Imposte Comunali Affini – s.r.l.'
$ws.Cells.Item(70, 2).NumberFormat = "@"
$ws.Cells.Item(70, 2).Value = '02478610583'

# Row 71
$ws.Cells.Item(71, 1).Value = 'Banco BPM Società per Azioni'
$ws.Cells.Item(71, 2).NumberFormat = "@"
$ws.Cells.Item(71, 2).Value = '09722490969'

# Row 72
$ws.Cells.Item(72, 1).Value = 'ARGO SOFTWARE SRL'
$ws.Cells.Item(72, 2).NumberFormat = "@"
$ws.Cells.Item(72, 2).Value = '00838520880'

# Row 73
$ws.Cells.Item(73, 1).Value = 'Softline srl'
$ws.Cells.Item(73, 2).NumberFormat = "@"
$ws.Cells.Item(73, 2).Value = '12299030150'

# Row 74
$ws.Cells.Item(74, 1).Value = 'San Marco SPA'
$ws.Cells.Item(74, 2).NumberFormat = "@"
$ws.Cells.Item(74, 2).Value = '04142440728'

# Row 75
$templateA.Copy()
$ws.Cells.Item(75, 1).PasteSpecial(-4122)
$templateB.Copy()
$ws.Cells.Item(75, 2).PasteSpecial(-4122)
$templateC.Copy()
$ws.Cells.Item(75, 3).PasteSpecial(-4122)
$ws.Cells.Item(75, 1).Value = 'Noviservice srl'
$ws.Cells.Item(75, 2).NumberFormat = "@"
$ws.Cells.Item(75, 2).Value = '02789990922'
$ws.Cells.Item(75, 3).Value = 1.0

# Row 76
$templateA.Copy()
$ws.Cells.Item(76, 1).PasteSpecial(-4122)
$templateB.Copy()
$ws.Cells.Item(76, 2).PasteSpecial(-4122)
$templateC.Copy()
$ws.Cells.Item(76, 3).PasteSpecial(-4122)
$ws.Cells.Item(76, 1).Value = 'Agenzia Italiana del Farmaco - AIFA'
$ws.Cells.Item(76, 2).NumberFormat = "@"
$ws.Cells.Item(76, 2).Value = '97345810580'
$ws.Cells.Item(76, 3).Value = 1.0

# Row 77
$templateA.Copy()
$ws.Cells.Item(77, 1).PasteSpecial(-4122)
$templateB.Copy()
$ws.Cells.Item(77, 2).PasteSpecial(-4122)
$templateC.Copy()
$ws.Cells.Item(77, 3).PasteSpecial(-4122)
$ws.Cells.Item(77, 1).Value = 'MegASP S.r.l.'
$ws.Cells.Item(77, 2).NumberFormat = "@"
$ws.Cells.Item(77, 2).Value = '09898030151'
$ws.Cells.Item(77, 3).Value = 1.0

$excel.CutCopyMode = 0
Write-Host "Edit complete"